# Update master data with new DAF values
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artikel FGR+")

$ws.Range("C3").Value = 0.6
$ws.Range("C4").Value = 0.5
$ws.Range("C5").Value = 0.7
$ws.Range("C16").Value = 0.7
$ws.Range("C17").Value = 0.7
$ws.Range("C18").Value = 1.1
$ws.Range("C19").Value = 0.8
$ws.Range("C20").Value = 0.9
$ws.Range("C23").Value = 0.9
$ws.Range("C24").Value = 1
$ws.Range("C25").Value = 0.9
$ws.Range("C27").Value = 0.8
$ws.Range("C28").Value = 1.3
$ws.Range("C29").Value = 1.1
$ws.Range("C30").Value = 0.7
$ws.Range("C31").Value = 1
$ws.Range("C33").Value = 2.2

# Update the view: scroll back to top-left and move selection to C34
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("C34").Select()
